$wb = $excel.ActiveWorkbook

# Update "Лист1" sheet
$ws1 = $wb.Worksheets.Item("Лист1")
$ws1.Range("B1").Value = "Не найдено"
$ws1.Range("B3").Value = "Не найдено"

# Update "BonpetData" sheet
$ws2 = $wb.Worksheets.Item("BonpetData")
$ws2.Range("B1").Value = "Не найдено"
$ws2.Range("B3").Value = "Не найдено"
